$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Days")

# Fill in rows 320 to 328 (A = sequential counter, B = date serial)
# continuing from row 319 (A319=318, B319=43418)
$startRow = 320
$endRow = 328
$startA = 319
$startB = 43419

for ($i = 0; $i -le ($endRow - $startRow); $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $startA + $i
    $ws.Cells.Item($row, 2).Value = $startB + $i

    # Copy style from the row above to preserve formatting
    $ws.Cells.Item($row - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($row - 1, 2).Copy() | Out-Null
    $ws.Cells.Item($row, 2).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}

# Update selection to match target state
$ws.Range("C324").Select() | Out-Null
